$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '73.206.14'
$ws.Range("E2").Value = '  +1.44%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '4.051.03'
$ws.Range("E3").Value = '  +0.56%  '

$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.47'
$ws.Range("E5").Value = '  +11.48%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.24'
$ws.Range("E6").Value = '  +2.03%  '

$ws.Range("E7").Value = '  -1.88%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.758'
$ws.Range("E9").Value = '  +0.77%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.171'
$ws.Range("E10").Value = '  -0.47%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.86'
$ws.Range("E11").Value = '  +12.63%  '

$ws.Range("E12").Value = '  -0.78%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.06'
$ws.Range("E13").Value = '  +3.69%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.704.08'
$ws.Range("E14").Value = '  +0.32%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.056.91'
$ws.Range("E15").Value = '  +0.52%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.30'
$ws.Range("E16").Value = '  +1.45%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.25'
$ws.Range("E17").Value = '  +4.65%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.76'
$ws.Range("E18").Value = '  +0.89%  '

$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.132'
$ws.Range("E19").Value = '  -0.64%  '

$ws.Range("B20").Value = 'WrappedBTC'
$ws.Range("C20").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '73.099.02'
$ws.Range("E20").Value = '  +1.28%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '443.06'
$ws.Range("E21").Value = '  +2.71%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.74'
$ws.Range("E22").Value = '  +11.94%  '

$ws.Range("E23").Value = '  -0.71%  '

$ws.Range("E24").Value = '  +1.84%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '14.47'
$ws.Range("E25").Value = '  +2.11%  '

$ws.Range("E26").Value = '  +21.51%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.46'
$ws.Range("E27").Value = '  +3.27%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.86'
$ws.Range("E28").Value = '  +1.32%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.97'
$ws.Range("E29").Value = '  +2.44%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '37.08'
$ws.Range("E30").Value = '  +0.82%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.97'
$ws.Range("E31").Value = '  +13.22%  '

$ws.Range("E32").Value = '  +4.30%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '13.67'
$ws.Range("E33").Value = '  +2.40%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '690.25'
$ws.Range("E34").Value = '  +2.08%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '49.05'
$ws.Range("E35").Value = '  +9.98%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '70.62'
$ws.Range("E36").Value = '  +7.04%  '

$ws.Range("E37").Value = '  +0.97%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0874'
$ws.Range("E38").Value = '  +5.38%  '

$ws.Range("E39").Value = '  -2.72%  '

$ws.Range("B40").Value = 'WEMIXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.40'
$ws.Range("E40").Value = '  +6.55%  '

$ws.Range("B41").Value = 'THORChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.33'
$ws.Range("E41").Value = '  +16.20%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.39'
$ws.Range("E42").Value = '  +0.35%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.998'
$ws.Range("E43").Value = '  -0.12%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0496'
$ws.Range("E44").Value = '  +2.29%  '

$ws.Range("E45").Value = '  +0.13%  '

$ws.Range("E46").Value = '  +0.61%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.71'
$ws.Range("E47").Value = '  +3.50%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.42'
$ws.Range("E48").Value = '  +0.08%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.53'
$ws.Range("E49").Value = '  +7.36%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.05'
$ws.Range("E50").Value = '  +1.36%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.18'
$ws.Range("E51").Value = '  +9.01%  '
